$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new measurement run was manually split into two rows: the old row 4
# (run #3, spanning 2018-12-18 19:30 .. 2018-12-19 04:45 UTC) actually
# contained two separate background periods, so it is split into a shorter
# run (new row 4) and the remainder (which becomes the old run #3, now
# pushed down into new row 5). Every run below it shifts down by one sheet
# row. Column A (the Run_number index) is regenerated/untouched - it still
# reads 1,2,3,... on each row - only the B:G measurement columns move down.
# ---------------------------------------------------------------------------

$lastOldRow = 65
$newRow4B = 43213.515625
$newRow4C = 43213.792361111111
$newRow4D = 0

# Shift columns B:F (value AND formatting, e.g. the yellow "weird data"
# highlight) down by one row, working bottom-up (row 65 -> 66, ..., row
# 4 -> 5) so a source row is always copied before it gets overwritten.
# The destination is cleared first: Range.Copy leaves a destination cell
# untouched wherever the source range has no cell at all (e.g. a row with
# no Comment), so without the pre-clear a stale value/style could survive.
for ($r = $lastOldRow; $r -ge 4; $r--) {
    $dst = $r + 1
    $ws.Range("B" + $dst + ":F" + $dst).ClearContents()
    $ws.Range("B" + $r + ":F" + $r).Copy($ws.Range("B" + $dst + ":F" + $dst))
}

# New row 4 holds the first (shorter) half of the split run; it carries no
# comment/sample-count, so clear whatever B4:F4 held after the copy above.
$ws.Range("B4:F4").ClearContents()
$ws.Cells.Item(4, 2).Value2 = $newRow4B
$ws.Cells.Item(4, 3).Value2 = $newRow4C
$ws.Cells.Item(4, 4).Value2 = $newRow4D

# Row 3's end time was corrected (the run actually ended earlier).
$ws.Cells.Item(3, 3).Value2 = 43213.466666666667

# The new last row (66) continues the Run_number sequence in column A.
$ws.Cells.Item($lastOldRow + 1, 1).Value2 = $lastOldRow

# Move the active selection to reflect where the edit was made.
$ws.Range("E4").Select()
